$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: insert a new blank column before N, shifting
# the old N/O/P (Late / Outstanding heading / Outstanding) columns right to
# O/P/Q, then restore the inserted column's width.
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = 10.2

# "Edit Repayment Schedule" sheet: move the selection to C15.
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Activate()
$wsEdit.Range("C15").Select()

# Make "Repayment schedule" the active/selected sheet & cell, which also
# updates the workbook's activeTab and clears tabSelected from whichever
# sheet previously had it (NewLoanInput).
$wsSchedule.Activate()
$wsSchedule.Range("K17").Select()
